$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Update the Approved/Rejected + ReasonToReject columns for two test case rows
$ws.Range("I8").Value = "Rejected"
$ws.Range("J8").Value = "Test step failed"

$ws.Range("I12").Value = "Rejected"
$ws.Range("J12").Value = "sds"

# Move the active selection to K16, matching the saved cursor position
$ws.Range("K16").Select()
